$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the title row (A1:C1): "4.c.1.1" -> "4.c.1" ---
$ws.Range("A1").Value = "4.с.1 Билим берүү мекемелерде диплом берилгем мугалимдердин үлүшү"
$ws.Range("B1").Value = "4.c.1 Доля дипломированных учителей в образовательных учереждениях"
$ws.Range("C1").Value = "4.c.1 Proportion of certified teachers in educational institutions"

# --- 2. Insert two new year columns (2011, 2012) right after column D (2010) ---
$ws.Columns("E:F").Insert()

$ws.Range("E4").Value = 2011
$ws.Range("F4").Value = 2012

$ws.Range("E5").Value = 89.6
$ws.Range("F5").Value = 87.5

$ws.Range("E6").Value = 93.3
$ws.Range("F6").Value = 93.9

$ws.Range("E7").Value = 92.8
$ws.Range("F7").Value = 94.1

# --- 3. Insert a new year column (2022) after the last existing year column (O / 2021) ---
$ws.Columns("P:P").Insert()

$ws.Range("P4").Value = 2022
$ws.Range("P5").Value = 94.2
$ws.Range("P6").Value = 96
$ws.Range("P6").NumberFormat = "0.0"
$ws.Range("P7").Value = 97.5

# --- 4. Update the selection on the sheet view ---
$ws.Range("Q4").Select()
